# Fix bug: "add anonymous transaction" always false + finish sending the
# donation detail workbook to the organisator.
#
# The previous, incomplete export only had two donation rows (a wrongly
# named "Michele Collet" donor on two lines). The fixed export lists the
# real transactions for the period: one anonymous donation (no name / no
# date recorded - the bug being fixed) followed by all of "Aurore Remy"'s
# individual donations, each with its own timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the formatting of the data block (currently rows 5:6) down to
# the new last row (18) before we touch the values, so every new row picks
# up the same thin-border look as the existing data rows. ---
$fmtSrc = $ws.Range("A5:C6")
$fmtSrc.Copy()
$ws.Range("A5:C18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Donation rows ---
# Row 5: anonymous donation (no donor name shown, no transaction date captured)
$ws.Range("A5").Value = "Anonymous donation"
$ws.Range("B5").Value = 2268
$ws.Range("C5").Value = ""

# Rows 6-18: Aurore Remy's individual donations
$donor = "Aurore Remy"

$ws.Range("A6").Value = $donor
$ws.Range("B6").Value = 45
$ws.Range("C6").Value = "2020-06-24 13:02:24"

$ws.Range("A7").Value = $donor
$ws.Range("B7").Value = 12
$ws.Range("C7").Value = "2020-06-24 13:09:04"

$ws.Range("A8").Value = $donor
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "2020-06-24 13:11:04"

$ws.Range("A9").Value = $donor
$ws.Range("B9").Value = 50
$ws.Range("C9").Value = "2020-06-24 13:15:57"

$ws.Range("A10").Value = $donor
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "2020-06-24 13:17:22"

$ws.Range("A11").Value = $donor
$ws.Range("B11").Value = 56
$ws.Range("C11").Value = "2020-06-24 14:33:02"

$ws.Range("A12").Value = $donor
$ws.Range("B12").Value = 123
$ws.Range("C12").Value = "2020-06-24 14:34:18"

$ws.Range("A13").Value = $donor
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "2020-06-24 14:35:07"

$ws.Range("A14").Value = $donor
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = "2020-06-24 14:36:57"

$ws.Range("A15").Value = $donor
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "2020-06-24 14:38:14"

$ws.Range("A16").Value = $donor
$ws.Range("B16").Value = 11
$ws.Range("C16").Value = "2020-06-24 21:14:52"

$ws.Range("A17").Value = $donor
$ws.Range("B17").Value = 6
$ws.Range("C17").Value = "2020-06-24 21:19:21"

$ws.Range("A18").Value = $donor
$ws.Range("B18").Value = 55
$ws.Range("C18").Value = "2020-06-24 21:20:00"

# --- Data rows now look left-aligned instead of centered ---
$ws.Range("A5:C18").HorizontalAlignment = -4131

# --- Keep the selection where Excel would leave it after entering the
# last cell of the table. ---
$ws.Range("C18").Select()
